$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1953.2
$ws.Range("I19").Value = 3703.4
$ws.Range("J19").Value = 640.55
$ws.Range("K19").Value = 3703.4
$ws.Range("L19").Value = 640.55
$ws.Range("M19").Value = -3528.4
$ws.Range("N19").Value = -990.55

$ws.Range("H55").Value = 88.07143000000001
$ws.Range("J55").Value = 96
$ws.Range("L55").Value = 96
$ws.Range("N55").Value = -524

$ws.Range("H98").Value = 1016.4737
$ws.Range("I98").Value = 1051.5625
$ws.Range("J98").Value = 829.3333
$ws.Range("K98").Value = 1051.5625
$ws.Range("L98").Value = 829.3333
$ws.Range("M98").Value = 446.4375
$ws.Range("N98").Value = -3825.3333

$ws.Range("H122").Value = 1016.4737
$ws.Range("I122").Value = 1051.5625
$ws.Range("J122").Value = 829.3333
$ws.Range("K122").Value = 3154.6875
$ws.Range("L122").Value = 2487.9999
$ws.Range("M122").Value = -704.6875
$ws.Range("N122").Value = -7387.9999

$ws.Range("H132").Value = 1935.2565
$ws.Range("I132").Value = 1328.2894
$ws.Range("K132").Value = 3984.8682
$ws.Range("M132").Value = -1454.8682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 10500
$ws.Range("J49").Value = 10500
$ws.Range("L49").Value = 10500
$ws.Range("N49").Value = -11020

$ws.Range("H63").Value = 1741.5
$ws.Range("I63").Value = 1672.8572
$ws.Range("K63").Value = 1672.8572
$ws.Range("M63").Value = -986.8571999999999

$ws.Range("H66").Value = 1741.5
$ws.Range("I66").Value = 1672.8572
$ws.Range("K66").Value = 8364.286
$ws.Range("M66").Value = -4932.286

$ws.Range("H74").Value = 194437.48
$ws.Range("I74").Value = 233316.81
$ws.Range("J74").Value = 75022.42999999999
$ws.Range("K74").Value = 233316.81
$ws.Range("L74").Value = 75022.42999999999
$ws.Range("M74").Value = -232442.81
$ws.Range("N74").Value = -76770.42999999999

$ws.Range("H77").Value = 194437.48
$ws.Range("I77").Value = 233316.81
$ws.Range("J77").Value = 75022.42999999999
$ws.Range("K77").Value = 1166584.05
$ws.Range("L77").Value = 375112.15
$ws.Range("M77").Value = -1162216.05
$ws.Range("N77").Value = -383848.15

$ws.Range("H102").Value = 11101.111
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 12113.75
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 12113.75
$ws.Range("M102").Value = -1378
$ws.Range("N102").Value = -15357.75

$ws.Range("H132").Value = 25221.61
$ws.Range("I132").Value = 40449.668
$ws.Range("J132").Value = 3581.7368
$ws.Range("K132").Value = 121349.004
$ws.Range("L132").Value = 10745.2104
$ws.Range("M132").Value = -118819.004
$ws.Range("N132").Value = -15805.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 873.8125
$ws.Range("I16").Value = 828.75
$ws.Range("J16").Value = 918.875
$ws.Range("K16").Value = 828.75
$ws.Range("L16").Value = 918.875
$ws.Range("M16").Value = -541.75
$ws.Range("N16").Value = -1492.875

$ws.Range("H99").Value = 64805.75
$ws.Range("J99").Value = 2555.4285
$ws.Range("L99").Value = 2555.4285
$ws.Range("N99").Value = -5551.4285

$ws.Range("H113").Value = 873.8125
$ws.Range("I113").Value = 828.75
$ws.Range("J113").Value = 918.875
$ws.Range("K113").Value = 828.75
$ws.Range("L113").Value = 918.875
$ws.Range("M113").Value = 1341.25
$ws.Range("N113").Value = -5258.875

$ws.Range("H126").Value = 64805.75
$ws.Range("J126").Value = 2555.4285
$ws.Range("L126").Value = 7666.2855
$ws.Range("N126").Value = -12606.2855

$ws.Range("H132").Value = 2278.44
$ws.Range("I132").Value = 1410.1875
$ws.Range("J132").Value = 3822
$ws.Range("K132").Value = 4230.5625
$ws.Range("L132").Value = 11466
$ws.Range("M132").Value = -1700.5625
$ws.Range("N132").Value = -16526

$ws.Range("H134").Value = 1641.84
$ws.Range("I134").Value = 1019.5294
$ws.Range("J134").Value = 2964.25
$ws.Range("K134").Value = 3058.5882
$ws.Range("L134").Value = 8892.75
$ws.Range("M134").Value = -523.5882000000001
$ws.Range("N134").Value = -13962.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 978.4545000000001
$ws.Range("I5").Value = 402.28
$ws.Range("J5").Value = 1736.579
$ws.Range("K5").Value = 1206.84
$ws.Range("L5").Value = 5209.737
$ws.Range("M5").Value = -1094.84
$ws.Range("N5").Value = -5433.737

$ws.Range("H12").Value = 37.588234
$ws.Range("I12").Value = 24.333334
$ws.Range("K12").Value = 73.00000199999999
$ws.Range("M12").Value = 99.99999800000001

$ws.Range("H86").Value = 686.6667
$ws.Range("I86").Value = 700
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 2100
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -914
$ws.Range("N86").Value = -3872

$ws.Range("H89").Value = 686.6667
$ws.Range("I89").Value = 700
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 6300
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = -372
$ws.Range("N89").Value = -16356

$ws.Range("H113").Value = 618.0333000000001
$ws.Range("I113").Value = 603
$ws.Range("J113").Value = 659.375
$ws.Range("K113").Value = 1809
$ws.Range("L113").Value = 1978.125
$ws.Range("M113").Value = 361
$ws.Range("N113").Value = -6318.125

$ws.Range("H116").Value = 4102.55
$ws.Range("I116").Value = 409.16666
$ws.Range("K116").Value = 1227.49998
$ws.Range("M116").Value = 2214.50002

$ws.Range("H135").Value = 978.4545000000001
$ws.Range("I135").Value = 402.28
$ws.Range("J135").Value = 1736.579
$ws.Range("K135").Value = 3620.52
$ws.Range("L135").Value = 15629.211
$ws.Range("M135").Value = -1085.52
$ws.Range("N135").Value = -20699.211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4178.5
$ws.Range("I80").Value = 4771.107
$ws.Range("K80").Value = 4771.107
$ws.Range("M80").Value = -3773.107

$ws.Range("H83").Value = 4178.5
$ws.Range("I83").Value = 4771.107
$ws.Range("K83").Value = 23855.535
$ws.Range("M83").Value = -18863.535

$ws.Range("H107").Value = 7519.5713
$ws.Range("I107").Value = 12862
$ws.Range("J107").Value = 396.33334
$ws.Range("K107").Value = 12862
$ws.Range("L107").Value = 396.33334
$ws.Range("M107").Value = -10942
$ws.Range("N107").Value = -4236.33334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 5150
$ws.Range("I18").Value = 5150
$ws.Range("K18").Value = 5150
$ws.Range("M18").Value = -4978

$ws.Range("H61").Value = 1757.7142
$ws.Range("I61").Value = 1460.8
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1460.8
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1258.8
$ws.Range("N61").Value = -2904

$ws.Range("H113").Value = 1757.7142
$ws.Range("I113").Value = 1460.8
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1460.8
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 709.2
$ws.Range("N113").Value = -6840
